# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (with
# value "stock" on every data row), inserted between "total" and "date" so
# the downstream JSON export can tag each property record with its
# category. Two company-name strings also get a stray embedded space
# removed as a side effect of the data regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H, shifting date/legislator_name/legislator_id
# (and their data) one column to the right.
$ws.Range("H1:H4").Insert(-4161)

# New header + values for the inserted "property_category" column.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"

# Fix stray embedded spaces in two company names.
$ws.Range("B3").Value = "久津實業公司（公開發行公司）"
$ws.Range("B4").Value = "博仁建設公司(公開發行公司）"
